# Fruta / hortaliza, semanal
# Inserts 3 new weekly price rows for "Vega Modelo de Temuco - Frutilla"
# above the existing row 258, shifting the previous rows 258-322 down to
# 261-325, and fills the 3 new rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 258..322 down by three rows, leaving three blank rows ready
# to receive the new weekly records.
$ws.Rows("258:260").Insert()

# --- New row 258 --------------------------------------------------------
$ws.Cells.Item(258, 1).Value = 10
$ws.Cells.Item(258, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(258, 3).Value = "La Araucanía"
$ws.Cells.Item(258, 4).Value = 45204
$ws.Cells.Item(258, 5).Value = 9
$ws.Cells.Item(258, 6).Value = "Fruta"
$ws.Cells.Item(258, 7).Value = 100101
$ws.Cells.Item(258, 8).Value = "Berries"
$ws.Cells.Item(258, 9).Value = 100112025
$ws.Cells.Item(258, 10).Value = "Frutilla"
$ws.Cells.Item(258, 11).Value = "Sin especificar"
$ws.Cells.Item(258, 12).Value = "Primera"
$ws.Cells.Item(258, 13).Value = 1100
$ws.Cells.Item(258, 14).Value = 14000
$ws.Cells.Item(258, 15).Value = 17000
$ws.Cells.Item(258, 16).Value = 15364
$ws.Cells.Item(258, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(258, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(258, 19).Value = 2195
$ws.Cells.Item(258, 20).Value = 7

# --- New row 259 --------------------------------------------------------
$ws.Cells.Item(259, 1).Value = 10
$ws.Cells.Item(259, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(259, 3).Value = "La Araucanía"
$ws.Cells.Item(259, 4).Value = 45204
$ws.Cells.Item(259, 5).Value = 9
$ws.Cells.Item(259, 6).Value = "Fruta"
$ws.Cells.Item(259, 7).Value = 100101
$ws.Cells.Item(259, 8).Value = "Berries"
$ws.Cells.Item(259, 9).Value = 100112025
$ws.Cells.Item(259, 10).Value = "Frutilla"
$ws.Cells.Item(259, 11).Value = "Sin especificar"
$ws.Cells.Item(259, 12).Value = "Segunda"
$ws.Cells.Item(259, 13).Value = 200
$ws.Cells.Item(259, 14).Value = 12000
$ws.Cells.Item(259, 15).Value = 12000
$ws.Cells.Item(259, 16).Value = 12000
$ws.Cells.Item(259, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(259, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(259, 19).Value = 1714
$ws.Cells.Item(259, 20).Value = 7

# --- New row 260 --------------------------------------------------------
$ws.Cells.Item(260, 1).Value = 10
$ws.Cells.Item(260, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value = "La Araucanía"
$ws.Cells.Item(260, 4).Value = 45204
$ws.Cells.Item(260, 5).Value = 9
$ws.Cells.Item(260, 6).Value = "Fruta"
$ws.Cells.Item(260, 7).Value = 100101
$ws.Cells.Item(260, 8).Value = "Berries"
$ws.Cells.Item(260, 9).Value = 100112025
$ws.Cells.Item(260, 10).Value = "Frutilla"
$ws.Cells.Item(260, 11).Value = "Sin especificar"
$ws.Cells.Item(260, 12).Value = "Tercera"
$ws.Cells.Item(260, 13).Value = 390
$ws.Cells.Item(260, 14).Value = 7000
$ws.Cells.Item(260, 15).Value = 8000
$ws.Cells.Item(260, 16).Value = 7359
$ws.Cells.Item(260, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(260, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(260, 19).Value = 1051
$ws.Cells.Item(260, 20).Value = 7

$ws.Range("A1").Select() | Out-Null
